$d = $word.ActiveDocument

function Find-ParagraphStartingWith($prefix) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $cand = $d.Paragraphs.Item($i)
        if ($cand.Range.Text -like "$prefix*") {
            return $cand
        }
    }
    return $null
}

$wdBrightGreen = 4

# --- 1) Re-highlight the "males/females" bullet from yellow to green -------
# Range.Font.HighlightColorIndex (as opposed to plain Range.HighlightColorIndex)
# stamps both the run's <w:rPr> and the trailing paragraph-mark's
# <w:pPr><w:rPr>, matching how the existing yellow highlight was authored
# (Word highlights the whole paragraph, mark included, when you select and
# highlight an entire bulleted line).
$pMales = Find-ParagraphStartingWith "There are also parts for males"
$pMales.Range.Font.HighlightColorIndex = $wdBrightGreen

# --- 2) Drop the stray "_GoBack" bookmark wherever Word last left it -------
# (there can only be one _GoBack at a time) so it ends up solely in its new
# spot below.
try {
    $stale = $d.Bookmarks.Item("_GoBack")
    if ($stale -ne $null) {
        $stale.Delete()
    }
} catch {
    # no pre-existing _GoBack -- nothing to remove
}

# --- 3) Re-add "_GoBack" to the now-empty paragraph right after the
#        "males/females" bullet (where Word's last edit actually landed). --
$pAfterMales = $pMales.Next()
$d.Bookmarks.Add("_GoBack", $pAfterMales.Range)

# --- 4) Highlight the "change costume" bullet green (it had no highlight
#        at all before). -----------------------------------------------------
$pCostume = Find-ParagraphStartingWith "Another constraint is that to allow people to change costume"
$pCostume.Range.Font.HighlightColorIndex = $wdBrightGreen
